# The "Prix Spot" sheet lists daily prices in columns B..FA (one column per
# day). A new day ("18-nov") needs to be inserted right before the block of
# "oct." columns, which currently starts at column DW. Inserting a whole
# column there shifts every existing column from DW..FA one slot to the
# right (DX..FB), matching the diff's dimension change A1:FA25 -> A1:FB25.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Prix Spot")

# Insert a new blank column before column DW (shifts DW:FA -> DX:FB).
$ws.Columns("DW").Insert()

# Header for the newly inserted column.
$ws.Range("DW1").Value = "18-nov"

# The new day has no data yet, so every hour row gets the placeholder "-"
# (same convention used by the other not-yet-available days in the sheet).
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 127).Value = "-"
}
